$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name and Link for rows whose ranking order changed
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("B39").Value = "Frax"
$ws.Range("C39").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"

# Update Price (D) and Volume(1h) (E) for all data rows.
# Price values that look like plain numbers must be forced to stay as
# text (matching the original inlineStr cells) by briefly switching the
# cell to a text number format, then restoring the default "Normal" style
# so no visible formatting change remains.
$ws.Range("D2").Value = "20.029.83"
$ws.Range("E2").Value = "  -7.53%  "
$ws.Range("D3").Value = "1.403.30"
$ws.Range("E3").Value = "  -8.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9933"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9979"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "271.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3670"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3105"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.01"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.003"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06481"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9925"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.381"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.127"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.26%  "
$ws.Range("D16").Value = "1.399.36"
$ws.Range("E16").Value = "  -9.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001015"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05697"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -13.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9969"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -15.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.563"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.244"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.91%  "
$ws.Range("D25").Value = "19.998.20"
$ws.Range("E25").Value = "  -7.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.232"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.09%  "
$ws.Range("D29").Value = "1.554.06"
$ws.Range("E29").Value = "  -9.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "109.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.101"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -15.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.285"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -12.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8252"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -12.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07645"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.356"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.431"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05772"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.793"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9964"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1900"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.079"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5284"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.496"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5111"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.765"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.036"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -11.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9979"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
